# new test case for OCI CU was added
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "Multiple Host" test-case row with a new "OCI Trial" one.
$ws.Range("A2").Value = "OCI Trial"

# Widen column A so the new (longer) label fits, and move the active
# selection onto the freshly-edited cell - mirrors what Excel does when a
# user types a new value into A2 and then widens the column.
$ws.Columns.Item(1).ColumnWidth = 20.729166
[void]$ws.Range("A2").Select()
